$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price cells whose new values look like plain numbers (e.g. "394.63") ---
# Force each one to Text format first so Excel stores the literal digit
# string instead of silently converting it to a floating point number.
# (NumberFormat/ClearFormats must be applied cell-by-cell: applying them to a
#  multi-area union Range here only affects the first area.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# --- Assign the numeric-looking Price values as literal text ---
$ws.Range("D5").Value = '394.63'
$ws.Range("D6").Value = '109.21'
$ws.Range("D8").Value = '0.567'
$ws.Range("D10").Value = '0.622'
$ws.Range("D11").Value = '39.30'
$ws.Range("D12").Value = '0.0972'
$ws.Range("D16").Value = '19.10'
$ws.Range("D21").Value = '3.32'
$ws.Range("D23").Value = '13.04'
$ws.Range("D24").Value = '301.35'
$ws.Range("D25").Value = '74.24'
$ws.Range("D27").Value = '28.05'
$ws.Range("D28").Value = '4.39'
$ws.Range("D34").Value = '11.01'
$ws.Range("D35").Value = '37.91'
$ws.Range("D36").Value = '0.0485'
$ws.Range("D37").Value = '2.12'
$ws.Range("D38").Value = '51.71'
$ws.Range("D40").Value = '0.999'
$ws.Range("D42").Value = '134.15'
$ws.Range("D44").Value = '17.22'
$ws.Range("D45").Value = '0.120'
$ws.Range("D46").Value = '3.97'
$ws.Range("D48").Value = '22.03'
$ws.Range("D50").Value = '2.07'
$ws.Range("D51").Value = '2.38'

# Restore default (General) formatting now that the text values are stored,
# so the cells keep looking like the rest of the sheet.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()

# --- Remaining cell updates (Coin names, Links, already non-numeric
#     Price strings, and Volume(1h) percentages) ---
$ws.Range("D2").Value = '57.205.41'
$ws.Range("E2").Value = '  +9.00%  '
$ws.Range("D3").Value = '3.257.67'
$ws.Range("E3").Value = '  +4.35%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("E6").Value = '  +4.88%  '
$ws.Range("D7").Value = '3.251.13'
$ws.Range("E7").Value = '  +4.16%  '
$ws.Range("E8").Value = '  +5.02%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("E10").Value = '  +3.40%  '
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("E12").Value = '  +12.86%  '
$ws.Range("E13").Value = '  +2.16%  '
$ws.Range("D14").Value = '3.773.87'
$ws.Range("E14").Value = '  +4.50%  '
$ws.Range("E15").Value = '  +4.14%  '
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").Value = '3.256.91'
$ws.Range("E17").Value = '  +5.42%  '
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("E19").Value = '  -6.16%  '
$ws.Range("D20").Value = '57.069.07'
$ws.Range("E20").Value = '  +9.20%  '
$ws.Range("E21").Value = '  +2.42%  '
$ws.Range("E22").Value = '  +8.29%  '
$ws.Range("E23").Value = '  +2.57%  '
$ws.Range("E24").Value = '  +11.98%  '
$ws.Range("E25").Value = '  +4.45%  '
$ws.Range("E26").Value = '  -2.79%  '
$ws.Range("E27").Value = '  +1.51%  '
$ws.Range("E28").Value = '  +3.72%  '
$ws.Range("E29").Value = '  -2.40%  '
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("E31").Value = '  -1.82%  '
$ws.Range("E32").Value = '  -0.03%  '
$ws.Range("E33").Value = '  +1.69%  '
$ws.Range("E34").Value = '  +0.84%  '
$ws.Range("E35").Value = '  +3.06%  '
$ws.Range("E36").Value = '  -1.31%  '
$ws.Range("E37").Value = '  +1.49%  '
$ws.Range("E38").Value = '  +3.56%  '
$ws.Range("E39").Value = '  +14.95%  '
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("E42").Value = '  +2.86%  '
$ws.Range("E43").Value = '  +1.89%  '
$ws.Range("B44").Value = 'Celestia'
$ws.Range("C44").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("E44").Value = '  +0.81%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("E45").Value = '  +3.11%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("E46").Value = '  -2.58%  '
$ws.Range("E47").Value = '  -3.10%  '
$ws.Range("E48").Value = '  -1.43%  '
$ws.Range("D49").Value = '2.154.74'
$ws.Range("E49").Value = '  +3.12%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("E51").Value = '  -2.87%  '
